$wb = $excel.ActiveWorkbook

# --- Sheet 1: Projections_NewlyAdded ---
$ws1 = $wb.Worksheets.Item("Projections_NewlyAdded")

$ws1.Range("E2").Value = "13932 (9487, 20613)"
$ws1.Range("F2").Value = "670 (532, 854)"
$ws1.Range("G2").Value = "105 (66, 163)"
$ws1.Range("H2").Value = "88 (50, 142)"
$ws1.Range("I2").Value = "17 (10, 25)"
$ws1.Range("J2").Value = "13 (8, 20)"
$ws1.Range("K2").Value = "4 (2, 8)"
$ws1.Range("L2").Value = "109 (70, 164)"
$ws1.Range("M2").Value = "92 (54, 145)"
$ws1.Range("N2").Value = "16 (9, 24)"
$ws1.Range("O2").Value = "13 (8, 19)"

$ws1.Range("E3").Value = "13881 (9405, 20619)"
$ws1.Range("F3").Value = "666 (524, 856)"
$ws1.Range("G3").Value = "104 (65, 162)"
$ws1.Range("H3").Value = "86 (49, 141)"
$ws1.Range("I3").Value = "16 (10, 25)"
$ws1.Range("J3").Value = "13 (8, 20)"
$ws1.Range("K3").Value = "4 (2, 8)"
$ws1.Range("L3").Value = "109 (70, 165)"
$ws1.Range("M3").Value = "92 (54, 145)"
$ws1.Range("N3").Value = "16 (10, 24)"
$ws1.Range("O3").Value = "13 (8, 19)"

$ws1.Range("E4").Value = "13821 (9310, 20695)"
$ws1.Range("F4").Value = "664 (518, 861)"
$ws1.Range("G4").Value = "102 (64, 161)"
$ws1.Range("H4").Value = "85 (49, 140)"
$ws1.Range("I4").Value = "16 (10, 25)"
$ws1.Range("J4").Value = "13 (8, 20)"
$ws1.Range("K4").Value = "4 (2, 8)"
$ws1.Range("L4").Value = "109 (69, 165)"
$ws1.Range("M4").Value = "92 (54, 145)"
$ws1.Range("N4").Value = "16 (10, 24)"
$ws1.Range("O4").Value = "13 (8, 20)"

$ws1.Range("E5").Value = "13768 (9204, 20783)"
$ws1.Range("F5").Value = "663 (512, 866)"
$ws1.Range("G5").Value = "101 (63, 161)"
$ws1.Range("H5").Value = "84 (48, 140)"
$ws1.Range("I5").Value = "16 (10, 25)"
$ws1.Range("J5").Value = "13 (8, 20)"
$ws1.Range("K5").Value = "4 (2, 8)"
$ws1.Range("L5").Value = "108 (69, 165)"
$ws1.Range("M5").Value = "92 (53, 145)"
$ws1.Range("N5").Value = "16 (10, 24)"
$ws1.Range("O5").Value = "13 (8, 20)"

$ws1.Range("E6").Value = "13742 (9101, 20852)"
$ws1.Range("F6").Value = "661 (507, 871)"
$ws1.Range("G6").Value = "100 (62, 160)"
$ws1.Range("H6").Value = "84 (47, 140)"
$ws1.Range("I6").Value = "16 (9, 24)"
$ws1.Range("J6").Value = "12 (7, 20)"
$ws1.Range("K6").Value = "4 (2, 8)"
$ws1.Range("L6").Value = "108 (69, 165)"
$ws1.Range("M6").Value = "91 (53, 145)"
$ws1.Range("N6").Value = "16 (10, 24)"
$ws1.Range("O6").Value = "13 (8, 20)"

$ws1.Range("E7").Value = "13716 (8995, 20921)"
$ws1.Range("F7").Value = "660 (501, 875)"
$ws1.Range("G7").Value = "100 (61, 160)"
$ws1.Range("H7").Value = "83 (47, 139)"
$ws1.Range("I7").Value = "15 (9, 24)"
$ws1.Range("J7").Value = "12 (7, 19)"
$ws1.Range("K7").Value = "4 (2, 8)"
$ws1.Range("L7").Value = "107 (68, 164)"
$ws1.Range("M7").Value = "91 (52, 144)"
$ws1.Range("N7").Value = "16 (10, 24)"
$ws1.Range("O7").Value = "13 (8, 20)"

$ws1.Range("E8").Value = "13698 (8886, 20984)"
$ws1.Range("F8").Value = "658 (495, 879)"
$ws1.Range("G8").Value = "99 (60, 160)"
$ws1.Range("H8").Value = "82 (46, 139)"
$ws1.Range("I8").Value = "15 (9, 24)"
$ws1.Range("J8").Value = "12 (7, 19)"
$ws1.Range("K8").Value = "4 (2, 8)"
$ws1.Range("L8").Value = "107 (67, 164)"
$ws1.Range("M8").Value = "90 (52, 144)"
$ws1.Range("N8").Value = "16 (10, 24)"
$ws1.Range("O8").Value = "13 (8, 20)"

$ws1.Range("E9").Value = "13697 (8808, 21068)"
$ws1.Range("F9").Value = "657 (489, 883)"
$ws1.Range("G9").Value = "98 (59, 160)"
$ws1.Range("H9").Value = "82 (46, 139)"
$ws1.Range("I9").Value = "15 (9, 24)"
$ws1.Range("J9").Value = "12 (7, 19)"
$ws1.Range("K9").Value = "4 (2, 8)"
$ws1.Range("L9").Value = "106 (66, 164)"
$ws1.Range("M9").Value = "89 (51, 143)"
$ws1.Range("N9").Value = "16 (10, 24)"
$ws1.Range("O9").Value = "13 (8, 20)"

$ws1.Range("E10").Value = "13708 (8738, 21189)"
$ws1.Range("F10").Value = "655 (483, 886)"
$ws1.Range("G10").Value = "98 (58, 160)"
$ws1.Range("H10").Value = "81 (45, 139)"
$ws1.Range("I10").Value = "15 (9, 24)"
$ws1.Range("J10").Value = "12 (7, 19)"
$ws1.Range("K10").Value = "4 (2, 8)"
$ws1.Range("L10").Value = "105 (66, 163)"
$ws1.Range("M10").Value = "88 (51, 143)"
$ws1.Range("N10").Value = "16 (9, 24)"
$ws1.Range("O10").Value = "13 (8, 20)"

$ws1.Range("E11").Value = "13712 (8663, 21262)"
$ws1.Range("F11").Value = "653 (477, 890)"
$ws1.Range("G11").Value = "97 (58, 160)"
$ws1.Range("H11").Value = "81 (44, 139)"
$ws1.Range("I11").Value = "15 (9, 24)"
$ws1.Range("J11").Value = "11 (7, 19)"
$ws1.Range("K11").Value = "4 (2, 8)"
$ws1.Range("L11").Value = "104 (65, 163)"
$ws1.Range("M11").Value = "88 (50, 142)"
$ws1.Range("N11").Value = "16 (9, 24)"
$ws1.Range("O11").Value = "13 (7, 20)"

$ws1.Range("E12").Value = "13706 (8595, 21333)"
$ws1.Range("F12").Value = "652 (472, 894)"
$ws1.Range("G12").Value = "96 (57, 160)"
$ws1.Range("H12").Value = "80 (44, 139)"
$ws1.Range("I12").Value = "14 (8, 24)"
$ws1.Range("J12").Value = "11 (7, 19)"
$ws1.Range("K12").Value = "4 (2, 8)"
$ws1.Range("L12").Value = "103 (64, 162)"
$ws1.Range("M12").Value = "87 (49, 142)"
$ws1.Range("N12").Value = "15 (9, 24)"
$ws1.Range("O12").Value = "12 (7, 20)"

$ws1.Range("E13").Value = "13700 (8522, 21372)"
$ws1.Range("F13").Value = "650 (466, 898)"
$ws1.Range("G13").Value = "96 (56, 160)"
$ws1.Range("H13").Value = "80 (43, 139)"
$ws1.Range("I13").Value = "14 (8, 23)"
$ws1.Range("J13").Value = "11 (6, 19)"
$ws1.Range("K13").Value = "4 (2, 8)"
$ws1.Range("L13").Value = "103 (63, 162)"
$ws1.Range("M13").Value = "86 (49, 142)"
$ws1.Range("N13").Value = "15 (9, 24)"
$ws1.Range("O13").Value = "12 (7, 19)"

$ws1.Range("E14").Value = "13675 (8445, 21401)"
$ws1.Range("F14").Value = "648 (461, 902)"
$ws1.Range("G14").Value = "95 (55, 160)"
$ws1.Range("H14").Value = "79 (43, 139)"
$ws1.Range("I14").Value = "14 (8, 23)"
$ws1.Range("J14").Value = "11 (6, 19)"
$ws1.Range("K14").Value = "4 (2, 8)"
$ws1.Range("L14").Value = "102 (62, 162)"
$ws1.Range("M14").Value = "85 (48, 141)"
$ws1.Range("N14").Value = "15 (9, 24)"
$ws1.Range("O14").Value = "12 (7, 19)"

$ws1.Range("E15").Value = "13662 (8372, 21441)"
$ws1.Range("F15").Value = "646 (455, 905)"
$ws1.Range("G15").Value = "95 (54, 160)"
$ws1.Range("H15").Value = "79 (42, 139)"
$ws1.Range("I15").Value = "14 (8, 23)"
$ws1.Range("J15").Value = "11 (6, 19)"
$ws1.Range("K15").Value = "4 (2, 8)"
$ws1.Range("L15").Value = "101 (62, 162)"
$ws1.Range("M15").Value = "85 (47, 141)"
$ws1.Range("N15").Value = "15 (9, 24)"
$ws1.Range("O15").Value = "12 (7, 19)"

$ws1.Range("E16").Value = "14406 (9850, 21212)"
$ws1.Range("F16").Value = "661 (531, 836)"
$ws1.Range("G16").Value = "114 (72, 176)"
$ws1.Range("H16").Value = "96 (55, 153)"
$ws1.Range("I16").Value = "18 (11, 27)"
$ws1.Range("J16").Value = "14 (8, 22)"
$ws1.Range("K16").Value = "4 (2, 7)"
$ws1.Range("L16").Value = "118 (75, 176)"
$ws1.Range("M16").Value = "100 (59, 155)"
$ws1.Range("N16").Value = "17 (10, 25)"
$ws1.Range("O16").Value = "13 (8, 21)"

$ws1.Range("E17").Value = "14295 (9725, 21151)"
$ws1.Range("F17").Value = "655 (521, 836)"
$ws1.Range("G17").Value = "113 (71, 175)"
$ws1.Range("H17").Value = "94 (54, 152)"
$ws1.Range("I17").Value = "17 (10, 27)"
$ws1.Range("J17").Value = "14 (8, 22)"
$ws1.Range("K17").Value = "4 (2, 7)"
$ws1.Range("L17").Value = "118 (75, 176)"
$ws1.Range("M17").Value = "100 (59, 155)"
$ws1.Range("N17").Value = "17 (10, 26)"
$ws1.Range("O17").Value = "13 (8, 21)"

$ws1.Range("E18").Value = "14181 (9588, 21146)"
$ws1.Range("F18").Value = "650 (512, 838)"
$ws1.Range("G18").Value = "111 (69, 174)"
$ws1.Range("H18").Value = "93 (53, 151)"
$ws1.Range("I18").Value = "17 (10, 27)"
$ws1.Range("J18").Value = "14 (8, 21)"
$ws1.Range("K18").Value = "4 (2, 7)"
$ws1.Range("L18").Value = "118 (75, 177)"
$ws1.Range("M18").Value = "100 (59, 155)"
$ws1.Range("N18").Value = "17 (10, 26)"
$ws1.Range("O18").Value = "14 (8, 21)"

$ws1.Range("E19").Value = "14080 (9447, 21154)"
$ws1.Range("F19").Value = "646 (504, 840)"
$ws1.Range("G19").Value = "110 (68, 173)"
$ws1.Range("H19").Value = "92 (52, 151)"
$ws1.Range("I19").Value = "17 (10, 27)"
$ws1.Range("J19").Value = "13 (8, 21)"
$ws1.Range("K19").Value = "4 (2, 7)"
$ws1.Range("L19").Value = "117 (75, 177)"
$ws1.Range("M19").Value = "99 (58, 155)"
$ws1.Range("N19").Value = "17 (10, 26)"
$ws1.Range("O19").Value = "14 (8, 21)"

$ws1.Range("E20").Value = "13983 (9298, 21146)"
$ws1.Range("F20").Value = "641 (495, 841)"
$ws1.Range("G20").Value = "109 (67, 173)"
$ws1.Range("H20").Value = "91 (51, 150)"
$ws1.Range("I20").Value = "17 (10, 26)"
$ws1.Range("J20").Value = "13 (8, 21)"
$ws1.Range("K20").Value = "4 (2, 7)"
$ws1.Range("L20").Value = "117 (74, 177)"
$ws1.Range("M20").Value = "99 (58, 155)"
$ws1.Range("N20").Value = "17 (10, 26)"
$ws1.Range("O20").Value = "14 (8, 21)"

$ws1.Range("E21").Value = "13903 (9153, 21142)"
$ws1.Range("F21").Value = "637 (487, 842)"
$ws1.Range("G21").Value = "107 (66, 172)"
$ws1.Range("H21").Value = "90 (51, 149)"
$ws1.Range("I21").Value = "16 (10, 26)"
$ws1.Range("J21").Value = "13 (8, 21)"
$ws1.Range("K21").Value = "4 (2, 7)"
$ws1.Range("L21").Value = "116 (74, 177)"
$ws1.Range("M21").Value = "98 (57, 155)"
$ws1.Range("N21").Value = "17 (10, 26)"
$ws1.Range("O21").Value = "14 (8, 21)"

$ws1.Range("E22").Value = "13829 (9021, 21135)"
$ws1.Range("F22").Value = "633 (479, 843)"
$ws1.Range("G22").Value = "106 (64, 171)"
$ws1.Range("H22").Value = "89 (50, 148)"
$ws1.Range("I22").Value = "16 (9, 26)"
$ws1.Range("J22").Value = "13 (7, 21)"
$ws1.Range("K22").Value = "4 (2, 7)"
$ws1.Range("L22").Value = "115 (73, 176)"
$ws1.Range("M22").Value = "97 (56, 155)"
$ws1.Range("N22").Value = "17 (10, 26)"
$ws1.Range("O22").Value = "14 (8, 21)"

$ws1.Range("E23").Value = "13755 (8891, 21106)"
$ws1.Range("F23").Value = "628 (470, 843)"
$ws1.Range("G23").Value = "105 (63, 170)"
$ws1.Range("H23").Value = "87 (49, 148)"
$ws1.Range("I23").Value = "16 (9, 26)"
$ws1.Range("J23").Value = "13 (7, 20)"
$ws1.Range("K23").Value = "4 (2, 7)"
$ws1.Range("L23").Value = "114 (72, 176)"
$ws1.Range("M23").Value = "96 (56, 154)"
$ws1.Range("N23").Value = "17 (10, 26)"
$ws1.Range("O23").Value = "13 (8, 21)"

$ws1.Range("E24").Value = "13678 (8757, 21087)"
$ws1.Range("F24").Value = "623 (463, 844)"
$ws1.Range("G24").Value = "103 (62, 169)"
$ws1.Range("H24").Value = "86 (48, 147)"
$ws1.Range("I24").Value = "16 (9, 25)"
$ws1.Range("J24").Value = "12 (7, 20)"
$ws1.Range("K24").Value = "4 (2, 7)"
$ws1.Range("L24").Value = "113 (71, 175)"
$ws1.Range("M24").Value = "95 (55, 154)"
$ws1.Range("N24").Value = "17 (10, 26)"
$ws1.Range("O24").Value = "13 (8, 21)"

$ws1.Range("E25").Value = "13616 (8629, 21045)"
$ws1.Range("F25").Value = "619 (455, 843)"
$ws1.Range("G25").Value = "102 (60, 169)"
$ws1.Range("H25").Value = "85 (47, 146)"
$ws1.Range("I25").Value = "15 (9, 25)"
$ws1.Range("J25").Value = "12 (7, 20)"
$ws1.Range("K25").Value = "4 (2, 7)"
$ws1.Range("L25").Value = "112 (70, 175)"
$ws1.Range("M25").Value = "95 (54, 153)"
$ws1.Range("N25").Value = "16 (10, 26)"
$ws1.Range("O25").Value = "13 (8, 21)"

$ws1.Range("E26").Value = "13532 (8505, 21011)"
$ws1.Range("F26").Value = "614 (447, 843)"
$ws1.Range("G26").Value = "101 (59, 168)"
$ws1.Range("H26").Value = "84 (46, 145)"
$ws1.Range("I26").Value = "15 (9, 25)"
$ws1.Range("J26").Value = "12 (7, 20)"
$ws1.Range("K26").Value = "4 (2, 7)"
$ws1.Range("L26").Value = "111 (69, 174)"
$ws1.Range("M26").Value = "94 (53, 152)"
$ws1.Range("N26").Value = "16 (10, 26)"
$ws1.Range("O26").Value = "13 (8, 21)"

$ws1.Range("E27").Value = "13451 (8383, 20951)"
$ws1.Range("F27").Value = "609 (439, 843)"
$ws1.Range("G27").Value = "100 (58, 167)"
$ws1.Range("H27").Value = "83 (45, 144)"
$ws1.Range("I27").Value = "15 (8, 25)"
$ws1.Range("J27").Value = "12 (7, 20)"
$ws1.Range("K27").Value = "4 (2, 7)"
$ws1.Range("L27").Value = "110 (68, 173)"
$ws1.Range("M27").Value = "93 (53, 152)"
$ws1.Range("N27").Value = "16 (9, 26)"
$ws1.Range("O27").Value = "13 (8, 21)"

$ws1.Range("E28").Value = "13355 (8265, 20888)"
$ws1.Range("F28").Value = "604 (431, 842)"
$ws1.Range("G28").Value = "98 (57, 166)"
$ws1.Range("H28").Value = "82 (44, 144)"
$ws1.Range("I28").Value = "15 (8, 25)"
$ws1.Range("J28").Value = "12 (6, 20)"
$ws1.Range("K28").Value = "4 (2, 7)"
$ws1.Range("L28").Value = "109 (67, 173)"
$ws1.Range("M28").Value = "92 (52, 151)"
$ws1.Range("N28").Value = "16 (9, 26)"
$ws1.Range("O28").Value = "13 (7, 21)"

$ws1.Range("E29").Value = "13276 (8150, 20806)"
$ws1.Range("F29").Value = "599 (424, 841)"
$ws1.Range("G29").Value = "97 (55, 164)"
$ws1.Range("H29").Value = "81 (42, 142)"
$ws1.Range("I29").Value = "14 (8, 25)"
$ws1.Range("J29").Value = "11 (6, 20)"
$ws1.Range("K29").Value = "4 (2, 7)"
$ws1.Range("L29").Value = "108 (66, 172)"
$ws1.Range("M29").Value = "91 (51, 150)"
$ws1.Range("N29").Value = "16 (9, 26)"
$ws1.Range("O29").Value = "13 (7, 21)"

# --- Sheet 2: Projections_HealthcareNeeds ---
$ws2 = $wb.Worksheets.Item("Projections_HealthcareNeeds")

$ws2.Range("E2").Value = "757 (387, 1408)"
$ws2.Range("F2").Value = "649 (284, 1281)"
$ws2.Range("G2").Value = "81 (35, 171)"
$ws2.Range("H2").Value = "76 (42, 132)"
$ws2.Range("I2").Value = "757 (387, 1408)"
$ws2.Range("J2").Value = "649 (284, 1281)"
$ws2.Range("K2").Value = "81 (35, 171)"
$ws2.Range("L2").Value = "76 (42, 132)"

$ws2.Range("E3").Value = "750 (382, 1404)"
$ws2.Range("F3").Value = "642 (278, 1274)"
$ws2.Range("G3").Value = "81 (35, 172)"
$ws2.Range("H3").Value = "77 (42, 133)"
$ws2.Range("I3").Value = "750 (382, 1404)"
$ws2.Range("J3").Value = "642 (278, 1274)"
$ws2.Range("K3").Value = "81 (35, 172)"
$ws2.Range("L3").Value = "77 (42, 133)"

$ws2.Range("E4").Value = "743 (376, 1398)"
$ws2.Range("F4").Value = "634 (273, 1267)"
$ws2.Range("G4").Value = "81 (34, 172)"
$ws2.Range("H4").Value = "77 (42, 133)"
$ws2.Range("I4").Value = "743 (376, 1398)"
$ws2.Range("J4").Value = "634 (273, 1267)"
$ws2.Range("K4").Value = "81 (34, 172)"
$ws2.Range("L4").Value = "77 (42, 133)"

$ws2.Range("E5").Value = "736 (370, 1391)"
$ws2.Range("F5").Value = "627 (267, 1260)"
$ws2.Range("G5").Value = "80 (34, 172)"
$ws2.Range("H5").Value = "76 (41, 133)"
$ws2.Range("I5").Value = "736 (370, 1391)"
$ws2.Range("J5").Value = "627 (267, 1260)"
$ws2.Range("K5").Value = "80 (34, 172)"
$ws2.Range("L5").Value = "76 (41, 133)"

$ws2.Range("E6").Value = "728 (364, 1385)"
$ws2.Range("F6").Value = "620 (263, 1253)"
$ws2.Range("G6").Value = "80 (34, 172)"
$ws2.Range("H6").Value = "76 (41, 133)"
$ws2.Range("I6").Value = "728 (364, 1385)"
$ws2.Range("J6").Value = "620 (263, 1253)"
$ws2.Range("K6").Value = "80 (34, 172)"
$ws2.Range("L6").Value = "76 (41, 133)"

$ws2.Range("E7").Value = "721 (359, 1378)"
$ws2.Range("F7").Value = "612 (258, 1246)"
$ws2.Range("G7").Value = "79 (33, 171)"
$ws2.Range("H7").Value = "75 (40, 133)"
$ws2.Range("I7").Value = "721 (359, 1378)"
$ws2.Range("J7").Value = "612 (258, 1246)"
$ws2.Range("K7").Value = "79 (33, 171)"
$ws2.Range("L7").Value = "75 (40, 133)"

$ws2.Range("E8").Value = "714 (353, 1371)"
$ws2.Range("F8").Value = "606 (254, 1240)"
$ws2.Range("G8").Value = "78 (33, 170)"
$ws2.Range("H8").Value = "75 (40, 132)"
$ws2.Range("I8").Value = "714 (353, 1371)"
$ws2.Range("J8").Value = "606 (254, 1240)"
$ws2.Range("K8").Value = "78 (33, 170)"
$ws2.Range("L8").Value = "75 (40, 132)"

$ws2.Range("E9").Value = "707 (348, 1366)"
$ws2.Range("F9").Value = "600 (249, 1233)"
$ws2.Range("G9").Value = "77 (32, 169)"
$ws2.Range("H9").Value = "74 (39, 131)"
$ws2.Range("I9").Value = "707 (348, 1366)"
$ws2.Range("J9").Value = "600 (249, 1233)"
$ws2.Range("K9").Value = "77 (32, 169)"
$ws2.Range("L9").Value = "74 (39, 131)"

$ws2.Range("E10").Value = "701 (343, 1359)"
$ws2.Range("F10").Value = "594 (245, 1227)"
$ws2.Range("G10").Value = "76 (32, 168)"
$ws2.Range("H10").Value = "73 (39, 130)"
$ws2.Range("I10").Value = "701 (343, 1359)"
$ws2.Range("J10").Value = "594 (245, 1227)"
$ws2.Range("K10").Value = "76 (32, 168)"
$ws2.Range("L10").Value = "73 (39, 130)"

$ws2.Range("E11").Value = "694 (338, 1352)"
$ws2.Range("F11").Value = "588 (241, 1222)"
$ws2.Range("G11").Value = "75 (31, 167)"
$ws2.Range("H11").Value = "72 (38, 129)"
$ws2.Range("I11").Value = "694 (338, 1352)"
$ws2.Range("J11").Value = "588 (241, 1222)"
$ws2.Range("K11").Value = "75 (31, 167)"
$ws2.Range("L11").Value = "72 (38, 129)"

$ws2.Range("E12").Value = "688 (333, 1347)"
$ws2.Range("F12").Value = "582 (237, 1216)"
$ws2.Range("G12").Value = "74 (30, 166)"
$ws2.Range("H12").Value = "71 (37, 129)"
$ws2.Range("I12").Value = "688 (333, 1347)"
$ws2.Range("J12").Value = "582 (237, 1216)"
$ws2.Range("K12").Value = "74 (30, 166)"
$ws2.Range("L12").Value = "71 (37, 129)"

$ws2.Range("E13").Value = "681 (328, 1341)"
$ws2.Range("F13").Value = "577 (234, 1212)"
$ws2.Range("G13").Value = "74 (30, 164)"
$ws2.Range("H13").Value = "70 (37, 128)"
$ws2.Range("I13").Value = "681 (328, 1341)"
$ws2.Range("J13").Value = "577 (234, 1212)"
$ws2.Range("K13").Value = "74 (30, 164)"
$ws2.Range("L13").Value = "70 (37, 128)"

$ws2.Range("E14").Value = "675 (322, 1337)"
$ws2.Range("F14").Value = "571 (230, 1207)"
$ws2.Range("G14").Value = "73 (29, 163)"
$ws2.Range("H14").Value = "69 (36, 127)"
$ws2.Range("I14").Value = "675 (322, 1337)"
$ws2.Range("J14").Value = "571 (230, 1207)"
$ws2.Range("K14").Value = "73 (29, 163)"
$ws2.Range("L14").Value = "69 (36, 127)"

$ws2.Range("E15").Value = "669 (317, 1332)"
$ws2.Range("F15").Value = "565 (226, 1202)"
$ws2.Range("G15").Value = "72 (29, 162)"
$ws2.Range("H15").Value = "69 (35, 126)"
$ws2.Range("I15").Value = "669 (317, 1332)"
$ws2.Range("J15").Value = "565 (226, 1202)"
$ws2.Range("K15").Value = "72 (29, 162)"
$ws2.Range("L15").Value = "69 (35, 126)"

$ws2.Range("E16").Value = "817 (419, 1520)"
$ws2.Range("F16").Value = "696 (312, 1381)"
$ws2.Range("G16").Value = "87 (37, 182)"
$ws2.Range("H16").Value = "81 (44, 142)"
$ws2.Range("I16").Value = "817 (419, 1520)"
$ws2.Range("J16").Value = "696 (312, 1381)"
$ws2.Range("K16").Value = "87 (37, 182)"
$ws2.Range("L16").Value = "81 (44, 142)"

$ws2.Range("E17").Value = "810 (413, 1517)"
$ws2.Range("F17").Value = "689 (306, 1377)"
$ws2.Range("G17").Value = "86 (37, 183)"
$ws2.Range("H17").Value = "82 (44, 143)"
$ws2.Range("I17").Value = "810 (413, 1517)"
$ws2.Range("J17").Value = "689 (306, 1377)"
$ws2.Range("K17").Value = "86 (37, 183)"
$ws2.Range("L17").Value = "82 (44, 143)"

$ws2.Range("E18").Value = "803 (408, 1511)"
$ws2.Range("F18").Value = "682 (301, 1371)"
$ws2.Range("G18").Value = "86 (37, 184)"
$ws2.Range("H18").Value = "81 (44, 143)"
$ws2.Range("I18").Value = "803 (408, 1511)"
$ws2.Range("J18").Value = "682 (301, 1371)"
$ws2.Range("K18").Value = "86 (37, 184)"
$ws2.Range("L18").Value = "81 (44, 143)"

$ws2.Range("E19").Value = "796 (402, 1505)"
$ws2.Range("F19").Value = "675 (296, 1364)"
$ws2.Range("G19").Value = "86 (36, 184)"
$ws2.Range("H19").Value = "81 (44, 143)"
$ws2.Range("I19").Value = "796 (402, 1505)"
$ws2.Range("J19").Value = "675 (296, 1364)"
$ws2.Range("K19").Value = "86 (36, 184)"
$ws2.Range("L19").Value = "81 (44, 143)"

$ws2.Range("E20").Value = "788 (396, 1497)"
$ws2.Range("F20").Value = "667 (290, 1357)"
$ws2.Range("G20").Value = "85 (36, 184)"
$ws2.Range("H20").Value = "81 (43, 143)"
$ws2.Range("I20").Value = "788 (396, 1497)"
$ws2.Range("J20").Value = "667 (290, 1357)"
$ws2.Range("K20").Value = "85 (36, 184)"
$ws2.Range("L20").Value = "81 (43, 143)"

$ws2.Range("E21").Value = "780 (390, 1490)"
$ws2.Range("F21").Value = "659 (285, 1349)"
$ws2.Range("G21").Value = "84 (35, 183)"
$ws2.Range("H21").Value = "80 (43, 142)"
$ws2.Range("I21").Value = "780 (390, 1490)"
$ws2.Range("J21").Value = "659 (285, 1349)"
$ws2.Range("K21").Value = "84 (35, 183)"
$ws2.Range("L21").Value = "80 (43, 142)"

$ws2.Range("E22").Value = "772 (383, 1482)"
$ws2.Range("F22").Value = "651 (280, 1341)"
$ws2.Range("G22").Value = "83 (35, 182)"
$ws2.Range("H22").Value = "79 (42, 142)"
$ws2.Range("I22").Value = "772 (383, 1482)"
$ws2.Range("J22").Value = "651 (280, 1341)"
$ws2.Range("K22").Value = "83 (35, 182)"
$ws2.Range("L22").Value = "79 (42, 142)"

$ws2.Range("E23").Value = "763 (376, 1473)"
$ws2.Range("F23").Value = "644 (274, 1332)"
$ws2.Range("G23").Value = "82 (34, 181)"
$ws2.Range("H23").Value = "78 (42, 141)"
$ws2.Range("I23").Value = "763 (376, 1473)"
$ws2.Range("J23").Value = "644 (274, 1332)"
$ws2.Range("K23").Value = "82 (34, 181)"
$ws2.Range("L23").Value = "78 (42, 141)"

$ws2.Range("E24").Value = "755 (370, 1464)"
$ws2.Range("F24").Value = "636 (268, 1323)"
$ws2.Range("G24").Value = "81 (33, 180)"
$ws2.Range("H24").Value = "77 (41, 140)"
$ws2.Range("I24").Value = "755 (370, 1464)"
$ws2.Range("J24").Value = "636 (268, 1323)"
$ws2.Range("K24").Value = "81 (33, 180)"
$ws2.Range("L24").Value = "77 (41, 140)"

$ws2.Range("E25").Value = "746 (362, 1455)"
$ws2.Range("F25").Value = "628 (263, 1313)"
$ws2.Range("G25").Value = "80 (33, 179)"
$ws2.Range("H25").Value = "76 (40, 139)"
$ws2.Range("I25").Value = "746 (362, 1455)"
$ws2.Range("J25").Value = "628 (263, 1313)"
$ws2.Range("K25").Value = "80 (33, 179)"
$ws2.Range("L25").Value = "76 (40, 139)"

$ws2.Range("E26").Value = "736 (355, 1445)"
$ws2.Range("F26").Value = "619 (257, 1304)"
$ws2.Range("G26").Value = "79 (32, 177)"
$ws2.Range("H26").Value = "75 (39, 138)"
$ws2.Range("I26").Value = "736 (355, 1445)"
$ws2.Range("J26").Value = "619 (257, 1304)"
$ws2.Range("K26").Value = "79 (32, 177)"
$ws2.Range("L26").Value = "75 (39, 138)"

$ws2.Range("E27").Value = "726 (348, 1436)"
$ws2.Range("F27").Value = "610 (251, 1295)"
$ws2.Range("G27").Value = "78 (31, 176)"
$ws2.Range("H27").Value = "74 (38, 137)"
$ws2.Range("I27").Value = "726 (348, 1436)"
$ws2.Range("J27").Value = "610 (251, 1295)"
$ws2.Range("K27").Value = "78 (31, 176)"
$ws2.Range("L27").Value = "74 (38, 137)"

$ws2.Range("E28").Value = "716 (341, 1425)"
$ws2.Range("F28").Value = "602 (245, 1284)"
$ws2.Range("G28").Value = "77 (31, 174)"
$ws2.Range("H28").Value = "73 (38, 136)"
$ws2.Range("I28").Value = "716 (341, 1425)"
$ws2.Range("J28").Value = "602 (245, 1284)"
$ws2.Range("K28").Value = "77 (31, 174)"
$ws2.Range("L28").Value = "73 (38, 136)"

$ws2.Range("E29").Value = "706 (334, 1413)"
$ws2.Range("F29").Value = "592 (240, 1273)"
$ws2.Range("G29").Value = "75 (30, 172)"
$ws2.Range("H29").Value = "72 (37, 134)"
$ws2.Range("I29").Value = "706 (334, 1413)"
$ws2.Range("J29").Value = "592 (240, 1273)"
$ws2.Range("K29").Value = "75 (30, 172)"
$ws2.Range("L29").Value = "72 (37, 134)"

